{"js": "// Replace each three-digit-division problem's text with its new value.\n// The document is a drill table of \"NNN\u00f7D=\" cells; every old problem\n// string is unique, so an exact-text search/replace per cell is safe.\nconst replacements = [\n  [\"254\u00f73=\", \"447\u00f78=\"],\n  [\"778\u00f77=\", \"555\u00f79=\"],\n  [\"809\u00f73=\", \"798\u00f77=\"],\n  [\"566\u00f72=\", \"583\u00f73=\"],\n  [\"828\u00f75=\", \"375\u00f73=\"],\n  [\"564\u00f79=\", \"222\u00f79=\"],\n  [\"408\u00f72=\", \"280\u00f73=\"],\n  [\"138\u00f76=\", \"944\u00f72=\"],\n  [\"694\u00f79=\", \"449\u00f74=\"],\n  [\"332\u00f78=\", \"732\u00f77=\"],\n  [\"561\u00f72=\", \"639\u00f76=\"],\n  [\"662\u00f79=\", \"588\u00f72=\"],\n  [\"114\u00f77=\", \"739\u00f74=\"],\n  [\"733\u00f76=\", \"741\u00f74=\"],\n  [\"307\u00f74=\", \"925\u00f76=\"],\n  [\"882\u00f72=\", \"724\u00f77=\"],\n  [\"401\u00f78=\", \"853\u00f77=\"],\n  [\"215\u00f76=\", \"574\u00f72=\"],\n  [\"117\u00f79=\", \"194\u00f78=\"],\n  [\"147\u00f73=\", \"512\u00f72=\"],\n  [\"281\u00f74=\", \"104\u00f76=\"],\n  [\"955\u00f73=\", \"891\u00f74=\"],\n  [\"902\u00f72=\", \"487\u00f78=\"],\n  [\"587\u00f77=\", \"500\u00f79=\"],\n  [\"726\u00f72=\", \"976\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-division problem's text with its new value.\n# The document is a drill table of \"NNN\u00f7D=\" cells; every old problem\n# string is unique, so an exact-text Find/Replace per cell is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"254\u00f73=\", \"447\u00f78=\"),\n    @(\"778\u00f77=\", \"555\u00f79=\"),\n    @(\"809\u00f73=\", \"798\u00f77=\"),\n    @(\"566\u00f72=\", \"583\u00f73=\"),\n    @(\"828\u00f75=\", \"375\u00f73=\"),\n    @(\"564\u00f79=\", \"222\u00f79=\"),\n    @(\"408\u00f72=\", \"280\u00f73=\"),\n    @(\"138\u00f76=\", \"944\u00f72=\"),\n    @(\"694\u00f79=\", \"449\u00f74=\"),\n    @(\"332\u00f78=\", \"732\u00f77=\"),\n    @(\"561\u00f72=\", \"639\u00f76=\"),\n    @(\"662\u00f79=\", \"588\u00f72=\"),\n    @(\"114\u00f77=\", \"739\u00f74=\"),\n    @(\"733\u00f76=\", \"741\u00f74=\"),\n    @(\"307\u00f74=\", \"925\u00f76=\"),\n    @(\"882\u00f72=\", \"724\u00f77=\"),\n    @(\"401\u00f78=\", \"853\u00f77=\"),\n    @(\"215\u00f76=\", \"574\u00f72=\"),\n    @(\"117\u00f79=\", \"194\u00f78=\"),\n    @(\"147\u00f73=\", \"512\u00f72=\"),\n    @(\"281\u00f74=\", \"104\u00f76=\"),\n    @(\"955\u00f73=\", \"891\u00f74=\"),\n    @(\"902\u00f72=\", \"487\u00f78=\"),\n    @(\"587\u00f77=\", \"500\u00f79=\"),\n    @(\"726\u00f72=\", \"976\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
